# Add team Wins/Losses/Ties record columns (AD:AF) to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, bordered, centered) used by the
# other header cells (e.g. A1) so the new headers look consistent.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Team record is the same for every player row (2-38): 89 wins, 73 losses, 1 tie.
$ws.Range("AD2:AD38").Value = 89
$ws.Range("AE2:AE38").Value = 73
$ws.Range("AF2:AF38").Value = 1
